$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22-25 down to 23-26
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with data (a new weekly entry)
$ws.Cells.Item(22, 1).Value = 5
$ws.Cells.Item(22, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(22, 3).Value = "Maule"
$ws.Cells.Item(22, 4).Value = 44476
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 7
$ws.Cells.Item(22, 6).Value = 300000000
$ws.Cells.Item(22, 7).Value = "Espárragos"
$ws.Cells.Item(22, 8).Value = "Verde"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 5000
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 1100
$ws.Cells.Item(22, 13).Value = 1040
$ws.Cells.Item(22, 14).Value = "`$/kilo"
$ws.Cells.Item(22, 15).Value = "Provincia de Linares"
$ws.Cells.Item(22, 16).Value = 1040
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = "Hortaliza"
